$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D=294; E=293.1000061035156; F=298.2000122070312; G=291.6000061035156; H=5894024},
    @{Row=3; D=294; E=293.1000061035156; F=298.2000122070312; G=291.6000061035156; H=5894024},
    @{Row=4; D=294; E=293.1000061035156; F=298.2000122070312; G=291.6000061035156; H=5894024},
    @{Row=5; D=296.7000122070312; E=289.7999877929688; F=297; G=288; H=5894024},
    @{Row=6; D=304.2000122070312; E=304.7999877929688; F=313.2000122070312; G=300.2999877929688; H=5894024},
    @{Row=7; D=298.5; E=297.8999938964844; F=318; G=296.3999938964844; H=5894024},
    @{Row=8; D=299.7000122070312; E=297.6000061035156; F=300; G=296.7000122070312; H=5894024},
    @{Row=9; D=450; E=392.1000061035156; F=450.6000061035156; G=303.2999877929688; H=5894024},
    @{Row=10; D=319.7999877929688; E=192.8999938964844; F=323.7000122070312; G=140.1000061035156; H=5894024},
    @{Row=11; D=150.6000061035156; E=147.6000061035156; F=151.8000030517578; G=119.4000015258789; H=5894024},
    @{Row=12; D=88.5; E=117; F=129.3000030517578; G=84.59999847412109; H=5894024},
    @{Row=13; D=102.3000030517578; E=116.0999984741211; F=117.5999984741211; G=91.5; H=5894024},
    @{Row=14; D=62.09999847412109; E=81.30000305175781; F=91.1999969482422; G=56.09999847412109; H=5894024},
    @{Row=15; D=68.09999847412109; E=67.5; F=73.19999694824219; G=60; H=5894024},
    @{Row=16; D=55.5; E=49.20000076293945; F=62.09999847412109; G=48.59999847412109; H=5894024},
    @{Row=17; D=17.15999984741211; E=12.39000034332275; F=17.70000076293945; G=9.600000381469728; H=5894024},
    @{Row=18; D=11.69999980926514; E=12.0600004196167; F=21.42000007629395; G=11.69999980926514; H=5894024},
    @{Row=19; D=9.899999618530272; E=8.430000305175781; F=10.19999980926514; G=7.050000190734863; H=5894024},
    @{Row=20; D=5.809999942779541; E=7.190000057220459; F=11.42000007629394; G=5.190000057220459; H=5894024},
    @{Row=21; D=5.539999961853027; E=7.909999847412109; F=9.54300022125244; G=5.349999904632568; H=5894024},
    @{Row=22; D=13.44999980926514; E=13.34000015258789; F=18.96999931335449; G=12.02999973297119; H=5894024},
    @{Row=23; D=6.610000133514404; E=8.260000228881836; F=10.80000019073486; G=5; H=5894024},
    @{Row=24; D=7.489999771118164; E=8.859999656677246; F=13.74600028991699; G=7.214000225067139; H=5894024}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = "RBOT"
}

